$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = " Proyecto ID"
$ws.Range("C1").Value = "Proyecto Descripción"

$ws.Range("E1").Select()
